# Updated cryptos list on Fri Jan 19 15:30:58 UTC 2024 with GitHub Actions
#
# Refreshes the Price / Volume(1h) columns (and, where two coins swapped
# rank, the Coin / Link columns too) to the latest scraped values.
#
# NOTE: every new value is written with a leading "'" (quote-prefix) so
# Excel stores it as literal text instead of auto-coercing number-shaped
# strings like "312.00" into the number 312. The Style reset afterwards
# drops the quote-prefix/text-format styling that assigning the value
# would otherwise leave behind, so the cell ends up plain text with the
# sheet's default (unstyled) formatting -- matching how the source data
# cells already looked.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "'41.481.87"
    "E2" = "'  -2.82%  "
    "D3" = "'2.464.48"
    "E3" = "'  -2.47%  "
    "E4" = "'  +0.93%  "
    "D5" = "'312.00"
    "E5" = "'  -0.85%  "
    "D6" = "'91.17"
    "E6" = "'  -7.44%  "
    "D7" = "'0.539"
    "E7" = "'  -4.12%  "
    "E8" = "'  +0.81%  "
    "D9" = "'0.489"
    "E9" = "'  -5.30%  "
    "D10" = "'32.58"
    "E10" = "'  -7.51%  "
    "D11" = "'0.0773"
    "E11" = "'  -3.43%  "
    "E12" = "'  -0.26%  "
    "D13" = "'2.858.37"
    "E13" = "'  -2.10%  "
    "D14" = "'6.79"
    "E14" = "'  -6.00%  "
    "B15" = "'WrappedEther"
    "C15" = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
    "D15" = "'2.495.63"
    "E15" = "'  -1.07%  "
    "B16" = "'Chainlink"
    "C16" = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
    "D16" = "'15.19"
    "E16" = "'  -0.15%  "
    "D17" = "'0.776"
    "E17" = "'  -4.26%  "
    "D18" = "'41.282.52"
    "E18" = "'  -3.30%  "
    "D19" = "'6.23"
    "E19" = "'  -5.36%  "
    "D20" = "'0.0₃0913"
    "E20" = "'  -2.96%  "
    "D21" = "'70.02"
    "E21" = "'  +1.47%  "
    "D22" = "'10.87"
    "E22" = "'  -10.51%  "
    "D23" = "'233.66"
    "E23" = "'  -3.40%  "
    "D24" = "'2.72"
    "E24" = "'  -4.77%  "
    "E25" = "'  +0.09%  "
    "D26" = "'1.86"
    "E26" = "'  -6.88%  "
    "D27" = "'23.76"
    "E27" = "'  -6.78%  "
    "E28" = "'  -0.66%  "
    "D29" = "'9.64"
    "E29" = "'  -3.61%  "
    "D30" = "'35.76"
    "E30" = "'  -4.62%  "
    "D31" = "'152.23"
    "E31" = "'  -1.95%  "
    "D32" = "'5.36"
    "E32" = "'  -9.17%  "
    "E33" = "'  -5.87%  "
    "E34" = "'  -3.41%  "
    "D35" = "'0.0746"
    "E35" = "'  -4.73%  "
    "D36" = "'17.35"
    "E36" = "'  -1.35%  "
    "D37" = "'2.96"
    "E37" = "'  -5.44%  "
    "D38" = "'1.82"
    "E38" = "'  -7.11%  "
    "D39" = "'0.112"
    "E39" = "'  -3.94%  "
    "D40" = "'0.0988"
    "E40" = "'  -8.44%  "
    "D41" = "'4.00"
    "E41" = "'  -5.46%  "
    "E42" = "'  +0.93%  "
    "D43" = "'19.25"
    "E43" = "'  -11.86%  "
    "D44" = "'1.952.78"
    "E44" = "'  -3.85%  "
    "D45" = "'0.0280"
    "E45" = "'  -5.44%  "
    "D46" = "'2.91"
    "E46" = "'  -9.39%  "
    "D47" = "'8.61"
    "E47" = "'  -2.65%  "
    "D48" = "'2.716.02"
    "E48" = "'  -1.97%  "
    "B49" = "'Aave"
    "C49" = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
    "D49" = "'94.83"
    "E49" = "'  -5.09%  "
    "B50" = "'ordi"
    "C50" = "'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
    "D50" = "'67.42"
    "E50" = "'  -6.11%  "
    "D51" = "'0.174"
    "E51" = "'  -7.65%  "
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}

Write-Host "Applied $($updates.Count) cell updates"
